$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column data types
$ws.Range("B1").Value = "varchar(50)"
$ws.Range("C1").Value = "varchar(50)"
$ws.Range("E1").Value = "varchar(50)"
$ws.Range("F1").Value = "int(255)"
$ws.Range("G1").Value = "int(6)"

# Row 2: column names
$ws.Range("A2").Value = "ID (입력x)"
$ws.Range("B2").Value = "법인코드"
$ws.Range("C2").Value = "bom코드"
$ws.Range("D2").Value = "계정코드"
$ws.Range("E2").Value = "version코드"
$ws.Range("F2").Value = "금액"
$ws.Range("G2").Value = "년월 ex) 200001"
